$d = $word.ActiveDocument

# The purchase order table: SR.NO | DESCRIPTION | UNIT OF MEASUREMENT | QTY | UNIT PRICE | AMOUNT
$tbl = $d.Tables.Item(1)

# Header row: the "UNIT OF MEASUREMENT" and "QTY" column headers were swapped.
$tbl.Cell(1, 3).Range.Text = "QTY"
$tbl.Cell(1, 4).Range.Text = "UNIT OF MEASUREMENT"

# Data row: the unit price / amount values now carry an "AED " currency prefix.
$tbl.Cell(2, 5).Range.Text = "AED 4,583.33"
$tbl.Cell(2, 6).Range.Text = "AED 55,000.00"
